# "Bundle make_ascii() to rebel()"
#
# Net effect (per the target OOXML diff): a new worksheet named "Sheet2" is
# inserted before the existing "Sheet1", carrying a duplicate of Sheet1's
# data (same dimension A1:G6, same values -> same shared-string usage
# doubles from 42 to 84 while uniqueCount stays 42) and keeping Sheet1's
# old selection (L2). The original "Sheet1" tab keeps its data, stays the
# active/selected sheet, and its selection moves to G26. The book window
# position also shifts (best-effort; not all hosts expose this).

$wb = $excel.ActiveWorkbook

# The sheet we are duplicating.
$source = $wb.Worksheets.Item("Sheet1")

# Insert a brand-new sheet directly before "Sheet1" -> Excel auto-names it
# "Sheet2" and it becomes the active sheet.
$copy = $wb.Worksheets.Add($source)

# Re-resolve both sheets by name: once a new sheet has been inserted, stale
# worksheet references captured before the insert stop resolving values.
$source = $wb.Worksheets.Item("Sheet1")
$copy = $wb.Worksheets.Item("Sheet2")

# Duplicate every cell value from Sheet1 onto the new Sheet2.
for ($r = 1; $r -le 6; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $copy.Cells.Item($r, $c).Value = $source.Cells.Item($r, $c).Value()
    }
}

# Sheet2 keeps the selection Sheet1 used to have...
[void]$copy.Range("L2").Select()

# ...while Sheet1 stays the active tab and picks up the new selection.
[void]$source.Activate()
[void]$source.Range("G26").Select()

# Best-effort: nudge the saved window position (not all hosts persist this).
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 29480
    $win.Top = 840
} catch {
}
